$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 684.26666
$ws.Range("J9").Value = 1413.6
$ws.Range("L9").Value = 1413.6
$ws.Range("N9").Value = -1751.6
$ws.Range("H18").Value = 1704.6
$ws.Range("I18").Value = 397.58334
$ws.Range("K18").Value = 397.58334
$ws.Range("M18").Value = -113.58334
$ws.Range("H32").Value = 2629.077
$ws.Range("I32").Value = 2619.25
$ws.Range("J32").Value = 2633.4443
$ws.Range("K32").Value = 2619.25
$ws.Range("L32").Value = 2633.4443
$ws.Range("M32").Value = -2293.25
$ws.Range("N32").Value = -3285.4443
$ws.Range("H33").Value = 1447.6316
$ws.Range("I33").Value = 397.18182
$ws.Range("J33").Value = 2892
$ws.Range("K33").Value = 397.18182
$ws.Range("L33").Value = 2892
$ws.Range("M33").Value = -168.18182
$ws.Range("N33").Value = -3350
$ws.Range("H38").Value = 2904.95
$ws.Range("I38").Value = 97.666664
$ws.Range("J38").Value = 4108.0713
$ws.Range("K38").Value = 292.999992
$ws.Range("L38").Value = 12324.2139
$ws.Range("M38").Value = 79.00000799999998
$ws.Range("N38").Value = -13068.2139
$ws.Range("H39").Value = 1363.1875
$ws.Range("I39").Value = 515.1667
$ws.Range("K39").Value = 1545.5001
$ws.Range("M39").Value = -1249.5001
$ws.Range("H48").Value = 4016.5833
$ws.Range("I48").Value = 5333
$ws.Range("K48").Value = 15999
$ws.Range("M48").Value = -15707
$ws.Range("H56").Value = 4016.5833
$ws.Range("I56").Value = 5333
$ws.Range("K56").Value = 15999
$ws.Range("M56").Value = -15465
$ws.Range("H76").Value = 7350
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 7828.5713
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 7828.5713
$ws.Range("M76").Value = -3685
$ws.Range("N76").Value = -8458.5713
$ws.Range("H79").Value = 7350
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 7828.5713
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 7828.5713
$ws.Range("M79").Value = -2908
$ws.Range("N79").Value = -10012.5713
$ws.Range("H96").Value = 686.9
$ws.Range("I96").Value = 846.75
$ws.Range("J96").Value = 47.5
$ws.Range("K96").Value = 2540.25
$ws.Range("L96").Value = 142.5
$ws.Range("M96").Value = -1167.25
$ws.Range("N96").Value = -2888.5
$ws.Range("H116").Value = 14649.5
$ws.Range("I116").Value = 15834.875
$ws.Range("K116").Value = 15834.875
$ws.Range("M116").Value = -12392.875
$ws.Range("H137").Value = 2397.6
$ws.Range("I137").Value = 2260.6316
$ws.Range("K137").Value = 6781.8948
$ws.Range("M137").Value = -4231.8948
$ws.Range("H138").Value = 3417.1072
$ws.Range("I138").Value = 2995.1667
$ws.Range("J138").Value = 4176.6
$ws.Range("K138").Value = 8985.500100000001
$ws.Range("L138").Value = 12529.8
$ws.Range("M138").Value = -3845.500100000001
$ws.Range("N138").Value = -22809.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 38873
$ws.Range("J133").Value = 38873
$ws.Range("L133").Value = 38873
$ws.Range("N133").Value = -43933

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2834
$ws.Range("I134").Value = 1740.8
$ws.Range("K134").Value = 5222.4
$ws.Range("M134").Value = -2687.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4394.579
$ws.Range("J58").Value = 6856.7144
$ws.Range("L58").Value = 6856.7144
$ws.Range("N58").Value = -7262.7144
$ws.Range("H99").Value = 2506.2666
$ws.Range("I99").Value = 2499.6428
$ws.Range("K99").Value = 2499.6428
$ws.Range("M99").Value = -1001.6428
$ws.Range("H126").Value = 2506.2666
$ws.Range("I126").Value = 2499.6428
$ws.Range("K126").Value = 7498.928400000001
$ws.Range("M126").Value = -5028.928400000001
$ws.Range("H134").Value = 33336552
$ws.Range("I134").Value = 38463716
$ws.Range("K134").Value = 115391148
$ws.Range("M134").Value = -115388613
$ws.Range("H136").Value = 4394.579
$ws.Range("J136").Value = 6856.7144
$ws.Range("L136").Value = 20570.1432
$ws.Range("N136").Value = -25670.1432
$ws.Range("H139").Value = 119999
$ws.Range("J139").Value = 119999
$ws.Range("L139").Value = 119999
$ws.Range("N139").Value = -130279
$ws.Range("H140").Value = 112386.29
$ws.Range("J140").Value = 112386.29
$ws.Range("L140").Value = 112386.29
$ws.Range("N140").Value = -122746.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 683.3684
$ws.Range("I2").Value = 104.92308
$ws.Range("J2").Value = 1936.6666
$ws.Range("K2").Value = 629.5384799999999
$ws.Range("L2").Value = 11619.9996
$ws.Range("M2").Value = -516.5384799999999
$ws.Range("N2").Value = -11845.9996
$ws.Range("H5").Value = 1428
$ws.Range("J5").Value = 1065
$ws.Range("L5").Value = 3195
$ws.Range("N5").Value = -3419
$ws.Range("H63").Value = 15712.25
$ws.Range("I63").Value = 19616.334
$ws.Range("K63").Value = 58849.00199999999
$ws.Range("M63").Value = -58100.00199999999
$ws.Range("H64").Value = 6874.375
$ws.Range("I64").Value = 6999.3335
$ws.Range("J64").Value = 6499.5
$ws.Range("K64").Value = 20998.0005
$ws.Range("L64").Value = 19498.5
$ws.Range("M64").Value = -20728.0005
$ws.Range("N64").Value = -20038.5
$ws.Range("H66").Value = 15712.25
$ws.Range("I66").Value = 19616.334
$ws.Range("K66").Value = 176547.006
$ws.Range("M66").Value = -172803.006
$ws.Range("H67").Value = 6874.375
$ws.Range("I67").Value = 6999.3335
$ws.Range("J67").Value = 6499.5
$ws.Range("K67").Value = 20998.0005
$ws.Range("L67").Value = 19498.5
$ws.Range("M67").Value = -20062.0005
$ws.Range("N67").Value = -21370.5
$ws.Range("H70").Value = 12430.071
$ws.Range("I70").Value = 8335.666999999999
$ws.Range("K70").Value = 25007.001
$ws.Range("M70").Value = -24692.001
$ws.Range("H73").Value = 12430.071
$ws.Range("I73").Value = 8335.666999999999
$ws.Range("K73").Value = 25007.001
$ws.Range("M73").Value = -23915.001
$ws.Range("H75").Value = 2728.9473
$ws.Range("I75").Value = 1713.8334
$ws.Range("J75").Value = 3197.4614
$ws.Range("K75").Value = 5141.5002
$ws.Range("L75").Value = 9592.3842
$ws.Range("M75").Value = -4143.5002
$ws.Range("N75").Value = -11588.3842
$ws.Range("H78").Value = 2728.9473
$ws.Range("I78").Value = 1713.8334
$ws.Range("J78").Value = 3197.4614
$ws.Range("K78").Value = 15424.5006
$ws.Range("L78").Value = 28777.1526
$ws.Range("M78").Value = -10432.5006
$ws.Range("N78").Value = -38761.1526
$ws.Range("H103").Value = 675
$ws.Range("I103").Value = 250
$ws.Range("J103").Value = 760
$ws.Range("K103").Value = 750
$ws.Range("L103").Value = 2280
$ws.Range("M103").Value = 129
$ws.Range("N103").Value = -4038
$ws.Range("H135").Value = 1428
$ws.Range("J135").Value = 1065
$ws.Range("L135").Value = 9585
$ws.Range("N135").Value = -14655

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 50029.5
$ws.Range("I39").Value = 8059
$ws.Range("J39").Value = 92000
$ws.Range("K39").Value = 8059
$ws.Range("L39").Value = 92000
$ws.Range("M39").Value = -7599
$ws.Range("N39").Value = -92920
$ws.Range("H93").Value = 830.76
$ws.Range("I93").Value = 623.5909
$ws.Range("K93").Value = 623.5909
$ws.Range("M93").Value = 624.4091
$ws.Range("H140").Value = 63085.75
$ws.Range("J140").Value = 63085.75
$ws.Range("L140").Value = 63085.75
$ws.Range("N140").Value = -73445.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 19499
$ws.Range("I44").Value = 19499
$ws.Range("K44").Value = 19499
$ws.Range("M44").Value = -18945
$ws.Range("H74").Value = 22900
$ws.Range("I74").Value = 25000
$ws.Range("K74").Value = 25000
$ws.Range("M74").Value = -24064
$ws.Range("H77").Value = 22900
$ws.Range("I77").Value = 25000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70320
$ws.Range("H100").Value = 1330.625
$ws.Range("J100").Value = 1378
$ws.Range("L100").Value = 2756
$ws.Range("N100").Value = -3838
